$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (F3 / 20 / YİYECEK) below the existing table
$ws.Range("A19").Value = "F3"
$ws.Range("B19").Value = 20
$ws.Range("C19").Value = "YİYECEK"

# Update the selected cell on the sheet, matching the new active selection
$ws.Range("C19").Select()
